# Updates to power sector for curtailment and moving CHP out of flexible
# resources; updates to fuel balancing priorities.
#
# On the "BGDPbES" sheet, row 15 corresponds to "crude oil" (BAU Guaranteed
# Dispatch Percentage by Electricity Source). Change the guaranteed dispatch
# percentage for every year (columns B:AK) from 0 to 1, and leave the
# worksheet with that range selected (matching the saved selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

$rng = $ws.Range("B15:AK15")
$rng.Value = 1

$rng.Select() | Out-Null
